{"js": "// Update anonymized contract outputs: renumber [[PERSON_n]] placeholders\n// in several paragraphs to match the target diff. Each affected paragraph\n// consists of a single run whose whole text is replaced.\n\nconst replacements = [\n  [\"[[PERSON_22]] \u2013 \u201epro [[PERSON_22]]\u201c\", \"[[PERSON_22]] \u2013 \u201epro [[PERSON_23]]\u201c\"],\n  [\"[[PERSON_23]] \u2013 \u201es [[PERSON_23]]\u201c\", \"[[PERSON_24]] \u2013 \u201es [[PERSON_24]]\u201c\"],\n  [\"[[PERSON_24]] \u2013 \u201ek [[PERSON_24]]\u201c\", \"[[PERSON_25]] \u2013 \u201ek [[PERSON_25]]\u201c\"],\n  [\"[[PERSON_25]] \u2013 \u201es [[PERSON_25]]\u201c\", \"[[PERSON_26]] \u2013 \u201es [[PERSON_26]]\u201c\"],\n  [\"[[PERSON_26]] \u2013 \u201eo [[PERSON_26]]\u201c\", \"[[PERSON_27]] \u2013 \u201eo [[PERSON_27]]\u201c\"],\n  [\"[[PERSON_27]] \u2013 \u201epro [[PERSON_27]]\u201c\", \"[[PERSON_28]] \u2013 \u201epro [[PERSON_28]]\u201c\"],\n  [\"[[PERSON_28]] \u2013 \u201es [[PERSON_28]]\u201c\", \"[[PERSON_29]] \u2013 \u201es [[PERSON_29]]\u201c\"],\n  [\"[[PERSON_29]] \u2013 \u201ek [[PERSON_29]]\u201c\", \"[[PERSON_30]] \u2013 \u201ek [[PERSON_30]]\u201c\"],\n  [\"[[PERSON_30]] \u2013 \u201es [[PERSON_30]]\u201c\", \"[[PERSON_31]] \u2013 \u201es [[PERSON_31]]\u201c\"],\n  [\"[[PERSON_31]] \u2013 \u201eo [[PERSON_31]]\u201c\", \"[[PERSON_32]] \u2013 \u201eo [[PERSON_32]]\u201c\"],\n  [\"[[PERSON_32]] \u2013 \u201eu [[PERSON_33]]\u201c\", \"[[PERSON_33]] \u2013 \u201eu [[PERSON_33]]\u201c\"],\n  [\"[[PERSON_37]] \u2013 \u201eo [[PERSON_38]]\u201c\", \"[[PERSON_37]] \u2013 \u201eo [[PERSON_37]]\u201c\"],\n  [\"[[PERSON_39]] \u2013 \u201es [[PERSON_40]]\u201c\", \"[[PERSON_38]] \u2013 \u201es [[PERSON_39]]\u201c\"],\n  [\"[[PERSON_41]] \u2013 \u201ek [[PERSON_41]]\u201c\", \"[[PERSON_40]] \u2013 \u201ek [[PERSON_40]]\u201c\"],\n  [\"[[PERSON_42]] \u2013 \u201eod [[PERSON_42]]\u201c\", \"[[PERSON_41]] \u2013 \u201eod [[PERSON_41]]\u201c\"],\n  [\"[[PERSON_43]] \u2013 \u201es [[PERSON_43]]\u201c\", \"[[PERSON_42]] \u2013 \u201es [[PERSON_42]]\u201c\"],\n  [\"[[PERSON_44]] \u2013 \u201eu [[PERSON_44]]\u201c\", \"[[PERSON_43]] \u2013 \u201eu [[PERSON_43]]\u201c\"],\n  [\"[[PERSON_45]] \u2013 \u201eo [[PERSON_45]]\u201c\", \"[[PERSON_44]] \u2013 \u201eo [[PERSON_44]]\u201c\"],\n  [\"[[PERSON_46]] \u2013 \u201ek [[PERSON_46]]\u201c\", \"[[PERSON_45]] \u2013 \u201ek [[PERSON_45]]\u201c\"],\n  [\n    \"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_2]], [[PERSON_5]], [[PERSON_26]] \u010di [[PERSON_47]].\",\n    \"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_2]], [[PERSON_5]], [[PERSON_27]] \u010di [[PERSON_46]].\",\n  ],\n  [\n    \"sv\u011bdek [[PERSON_30]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\",\n    \"sv\u011bdek [[PERSON_31]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\",\n  ],\n  [\n    \"tlumo\u010dn\u00edk [[PERSON_32]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\",\n    \"tlumo\u010dn\u00edk [[PERSON_33]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\",\n  ],\n  [\n    \"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_38]],\",\n    \"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_37]],\",\n  ],\n  [\n    \"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_31]].\",\n    \"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_32]].\",\n  ],\n  [\n    \"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_17]], [[PERSON_21]] a [[PERSON_45]].\",\n    \"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_17]], [[PERSON_21]] a [[PERSON_44]].\",\n  ],\n  [\n    \"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_46]]),\",\n    \"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_45]]),\",\n  ],\n  [\"[[PERSON_43]] (\u201ev\u00fdslech [[PERSON_43]]\u201c),\", \"[[PERSON_42]] (\u201ev\u00fdslech [[PERSON_42]]\u201c),\"],\n  [\"[[PERSON_39]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_39]]\u201c),\", \"[[PERSON_38]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_38]]\u201c),\"],\n  [\"PhDr. [[PERSON_31]] \u2013 psychologick\u00fd posudek,\", \"PhDr. [[PERSON_32]] \u2013 psychologick\u00fd posudek,\"],\n  [\"MUDr. [[PERSON_26]] \u2013 posudek z traumatologie,\", \"MUDr. [[PERSON_27]] \u2013 posudek z traumatologie,\"],\n  [\n    \"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_22]], [[PERSON_41]] nebo [[PERSON_28]].\",\n    \"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_22]], [[PERSON_40]] nebo [[PERSON_29]].\",\n  ],\n  [\"[[PERSON_42]],\", \"[[PERSON_41]],\"],\n  [\"[[PERSON_47]],\", \"[[PERSON_46]],\"],\n  [\"[[PERSON_25]],\", \"[[PERSON_26]],\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Build a lookup from the current (old) paragraph text to its replacement,\n// and apply each replacement exactly once (in document order) so that\n// paragraphs with identical old text - if any - are not double-replaced by\n// the same mapping entry and every mapping entry is consumed at most once.\nconst pending = replacements.map(([oldText, newText]) => ({\n  oldText,\n  newText,\n  done: false,\n}));\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  for (const entry of pending) {\n    if (!entry.done && text === entry.oldText) {\n      para.insertText(entry.newText, Word.InsertLocation.replace);\n      entry.done = true;\n      break;\n    }\n  }\n}\n\nawait context.sync();\n\nconst missed = pending.filter((e) => !e.done).map((e) => e.oldText);\nif (missed.length > 0) {\n  throw new Error(\"Unmatched replacements: \" + JSON.stringify(missed));\n}\n", "ps1": "# Update anonymized contract outputs: renumber [[PERSON_n]] placeholders\n# in several paragraphs to match the target diff. Each affected paragraph\n# consists of a single run whose whole text is replaced.\n\n$d = $word.ActiveDocument\n\n$map = @(\n    @(\"[[PERSON_22]] \u2013 \u201epro [[PERSON_22]]\u201c\", \"[[PERSON_22]] \u2013 \u201epro [[PERSON_23]]\u201c\"),\n    @(\"[[PERSON_23]] \u2013 \u201es [[PERSON_23]]\u201c\", \"[[PERSON_24]] \u2013 \u201es [[PERSON_24]]\u201c\"),\n    @(\"[[PERSON_24]] \u2013 \u201ek [[PERSON_24]]\u201c\", \"[[PERSON_25]] \u2013 \u201ek [[PERSON_25]]\u201c\"),\n    @(\"[[PERSON_25]] \u2013 \u201es [[PERSON_25]]\u201c\", \"[[PERSON_26]] \u2013 \u201es [[PERSON_26]]\u201c\"),\n    @(\"[[PERSON_26]] \u2013 \u201eo [[PERSON_26]]\u201c\", \"[[PERSON_27]] \u2013 \u201eo [[PERSON_27]]\u201c\"),\n    @(\"[[PERSON_27]] \u2013 \u201epro [[PERSON_27]]\u201c\", \"[[PERSON_28]] \u2013 \u201epro [[PERSON_28]]\u201c\"),\n    @(\"[[PERSON_28]] \u2013 \u201es [[PERSON_28]]\u201c\", \"[[PERSON_29]] \u2013 \u201es [[PERSON_29]]\u201c\"),\n    @(\"[[PERSON_29]] \u2013 \u201ek [[PERSON_29]]\u201c\", \"[[PERSON_30]] \u2013 \u201ek [[PERSON_30]]\u201c\"),\n    @(\"[[PERSON_30]] \u2013 \u201es [[PERSON_30]]\u201c\", \"[[PERSON_31]] \u2013 \u201es [[PERSON_31]]\u201c\"),\n    @(\"[[PERSON_31]] \u2013 \u201eo [[PERSON_31]]\u201c\", \"[[PERSON_32]] \u2013 \u201eo [[PERSON_32]]\u201c\"),\n    @(\"[[PERSON_32]] \u2013 \u201eu [[PERSON_33]]\u201c\", \"[[PERSON_33]] \u2013 \u201eu [[PERSON_33]]\u201c\"),\n    @(\"[[PERSON_37]] \u2013 \u201eo [[PERSON_38]]\u201c\", \"[[PERSON_37]] \u2013 \u201eo [[PERSON_37]]\u201c\"),\n    @(\"[[PERSON_39]] \u2013 \u201es [[PERSON_40]]\u201c\", \"[[PERSON_38]] \u2013 \u201es [[PERSON_39]]\u201c\"),\n    @(\"[[PERSON_41]] \u2013 \u201ek [[PERSON_41]]\u201c\", \"[[PERSON_40]] \u2013 \u201ek [[PERSON_40]]\u201c\"),\n    @(\"[[PERSON_42]] \u2013 \u201eod [[PERSON_42]]\u201c\", \"[[PERSON_41]] \u2013 \u201eod [[PERSON_41]]\u201c\"),\n    @(\"[[PERSON_43]] \u2013 \u201es [[PERSON_43]]\u201c\", \"[[PERSON_42]] \u2013 \u201es [[PERSON_42]]\u201c\"),\n    @(\"[[PERSON_44]] \u2013 \u201eu [[PERSON_44]]\u201c\", \"[[PERSON_43]] \u2013 \u201eu [[PERSON_43]]\u201c\"),\n    @(\"[[PERSON_45]] \u2013 \u201eo [[PERSON_45]]\u201c\", \"[[PERSON_44]] \u2013 \u201eo [[PERSON_44]]\u201c\"),\n    @(\"[[PERSON_46]] \u2013 \u201ek [[PERSON_46]]\u201c\", \"[[PERSON_45]] \u2013 \u201ek [[PERSON_45]]\u201c\"),\n    @(\"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_2]], [[PERSON_5]], [[PERSON_26]] \u010di [[PERSON_47]].\", \"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_2]], [[PERSON_5]], [[PERSON_27]] \u010di [[PERSON_46]].\"),\n    @(\"sv\u011bdek [[PERSON_30]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\", \"sv\u011bdek [[PERSON_31]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\"),\n    @(\"tlumo\u010dn\u00edk [[PERSON_32]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\", \"tlumo\u010dn\u00edk [[PERSON_33]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\"),\n    @(\"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_38]],\", \"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_37]],\"),\n    @(\"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_31]].\", \"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_32]].\"),\n    @(\"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_17]], [[PERSON_21]] a [[PERSON_45]].\", \"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_17]], [[PERSON_21]] a [[PERSON_44]].\"),\n    @(\"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_46]]),\", \"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_45]]),\"),\n    @(\"[[PERSON_43]] (\u201ev\u00fdslech [[PERSON_43]]\u201c),\", \"[[PERSON_42]] (\u201ev\u00fdslech [[PERSON_42]]\u201c),\"),\n    @(\"[[PERSON_39]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_39]]\u201c),\", \"[[PERSON_38]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_38]]\u201c),\"),\n    @(\"PhDr. [[PERSON_31]] \u2013 psychologick\u00fd posudek,\", \"PhDr. [[PERSON_32]] \u2013 psychologick\u00fd posudek,\"),\n    @(\"MUDr. [[PERSON_26]] \u2013 posudek z traumatologie,\", \"MUDr. [[PERSON_27]] \u2013 posudek z traumatologie,\"),\n    @(\"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_22]], [[PERSON_41]] nebo [[PERSON_28]].\", \"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_22]], [[PERSON_40]] nebo [[PERSON_29]].\"),\n    @(\"[[PERSON_42]],\", \"[[PERSON_41]],\"),\n    @(\"[[PERSON_47]],\", \"[[PERSON_46]],\"),\n    @(\"[[PERSON_25]],\", \"[[PERSON_26]],\")\n)\n\n$done = @()\nfor ($k = 0; $k -lt $map.Length; $k++) {\n    $done += $false\n}\n\n$n = $d.Paragraphs.Count\nfor ($i = 1; $i -le $n; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    # Exclude the trailing paragraph mark from both the comparison and the\n    # replacement so paragraph formatting (numbering, paraId, ...) survives.\n    [void]$r.MoveEnd(1, -1)\n    $text = $r.Text\n\n    for ($j = 0; $j -lt $map.Length; $j++) {\n        if (-not $done[$j]) {\n            $oldText = $map[$j][0]\n            if ($text -eq $oldText) {\n                $r.Text = $map[$j][1]\n                $done[$j] = $true\n                break\n            }\n        }\n    }\n}\n\n$missed = 0\nfor ($j = 0; $j -lt $map.Length; $j++) {\n    if (-not $done[$j]) {\n        $missed = $missed + 1\n        \"Missed: \" + $map[$j][0]\n    }\n}\n\"Replacements applied: \" + ($map.Length - $missed) + \" / \" + $map.Length\n"}
